$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 29 ("Oakleigh"), shifting it and the
# following "Sandringham Line" row down by one, then populate the new row
# with the "Nunawading" exposure site entry.
$ws.Rows.Item(29).Insert()

$ws.Range("A29").Value = "Nunawading"
$ws.Range("B29").Value = "Good Guys, 7/372 Whitehorse Road, Nunawading"
$ws.Range("C29").Value = "29/12/20, 9:30am-10:00am"
$ws.Range("D29").Value = "Case shopped in store"
